$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates (WORKFLOWINSTANCEID, IDDOCUMENTO, DATAEVENTO)
$ws.Range("D2").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.be7ba7c0c7^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721646186259"
$ws.Range("F2").Value = "22-07-2024:13:03:07"

# Row 3 updates (WORKFLOWINSTANCEID, IDDOCUMENTO, DATAEVENTO)
$ws.Range("D3").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.1a94c8365b^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E3").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721646174009"
$ws.Range("F3").Value = "22-07-2024:13:02:56"

$wb.Save()
